$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Court_Reports")

# Row 2
$ws.Range("B2").Value = "London"
$ws.Range("C2").Value = "ND01"
$ws.Range("E2").Value = "KNS"
$ws.Range("G2").Value = "WMT"
$ws.Range("H2").Value = "Court"
$ws.Range("I2").Value = "John"
$ws.Range("J2").Value = 2001
$ws.Range("K2").Value = "Z"
$ws.Range("L2").Value = "2001|WMT|Z"
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 3
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 7

# Row 3
$ws.Range("B3").Value = "London"
$ws.Range("C3").Value = "ND01"
$ws.Range("D3").Value = "CR LDU"
$ws.Range("E3").Value = "CRLDU"
$ws.Range("F3").Value = "CR Team"
$ws.Range("G3").Value = "CRP"
$ws.Range("H3").Value = "Report"
$ws.Range("I3").Value = "Simon"
$ws.Range("J3").Value = 2002
$ws.Range("K3").Value = "C"
$ws.Range("L3").Value = "2002|WMT|C"
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = 8

# Update the selection shown on the sheet view
$ws.Range("A2:XFD3").Select()
